$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T-H曲線")
$ws.Name = "T-H曲線_p0-0-0"
